$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so values like "580.88" or
# "66.938.68" are stored as text, matching the source data (inline strings),
# rather than being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.938.68"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").Value = "3.105.06"
$ws.Range("E3").Value = "  +5.24%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "580.88"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").Value = "173.15"
$ws.Range("E6").Value = "  +6.92%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.100.84"
$ws.Range("E8").Value = "  +5.12%  "
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").Value = "6.48"
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("D11").Value = "0.156"
$ws.Range("E11").Value = "  +4.15%  "
$ws.Range("E12").Value = "  +4.82%  "
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").Value = "37.24"
$ws.Range("E14").Value = "  +6.99%  "
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "3.619.84"
$ws.Range("E16").Value = "  +5.11%  "
$ws.Range("D17").Value = "66.910.77"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").Value = "3.109.40"
$ws.Range("E19").Value = "  +5.29%  "
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").Value = "485.32"
$ws.Range("E21").Value = "  +8.92%  "
$ws.Range("D22").Value = "0.714"
$ws.Range("E22").Value = "  +2.61%  "
$ws.Range("D23").Value = "7.52"
$ws.Range("E23").Value = "  +3.19%  "
$ws.Range("D24").Value = "84.06"
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  +6.48%  "
$ws.Range("D26").Value = "13.15"
$ws.Range("E26").Value = "  +7.43%  "
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "7.98"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").Value = "2.40"
$ws.Range("E30").Value = "  -4.68%  "
$ws.Range("D31").Value = "2.69"
$ws.Range("E31").Value = "  +3.67%  "
$ws.Range("D32").Value = "0.0000101"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").Value = "28.79"
$ws.Range("E33").Value = "  +6.13%  "
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "5.91"
$ws.Range("E36").Value = "  +3.51%  "
$ws.Range("D37").Value = "0.996"
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("E38").Value = "  +4.15%  "
$ws.Range("D39").Value = "2.12"
$ws.Range("E39").Value = "  +6.90%  "
$ws.Range("E40").Value = "  +4.94%  "
$ws.Range("D41").Value = "50.16"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "8.66"
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").Value = "0.0362"
$ws.Range("E45").Value = "  +3.20%  "
$ws.Range("D46").Value = "2.834.20"
$ws.Range("E46").Value = "  +5.78%  "
$ws.Range("D47").Value = "384.65"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").Value = "135.07"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "24.89"
$ws.Range("E50").Value = "  +4.40%  "
$ws.Range("E51").Value = "  +3.24%  "

# Restore the default cell style on column D so no stray number-format
# style is left behind (keeps formatting identical to the original file).
$ws.Range("D2:D51").Style = "Normal"
